$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Weekly update: shift Fecha/Volumen/Precio columns (D,J,K,L,M,P) down by one row
# for rows 63-108, and append new row 109 with the previous row108 data.

$ws.Range("D63").Value = 44447
$ws.Range("J63").Value = 3400
$ws.Range("K63").Value = 900
$ws.Range("L63").Value = 1000
$ws.Range("M63").Value = 950
$ws.Range("P63").Value = 158

$ws.Range("D64").Value = 44267
$ws.Range("J64").Value = 2400
$ws.Range("K64").Value = 900
$ws.Range("L64").Value = 1000
$ws.Range("M64").Value = 950
$ws.Range("P64").Value = 158

$ws.Range("D65").Value = 44300
$ws.Range("J65").Value = 3200
$ws.Range("K65").Value = 900
$ws.Range("L65").Value = 1000
$ws.Range("M65").Value = 950
$ws.Range("P65").Value = 158

$ws.Range("D66").Value = 44277
$ws.Range("J66").Value = 2400
$ws.Range("K66").Value = 900
$ws.Range("L66").Value = 1000
$ws.Range("M66").Value = 950
$ws.Range("P66").Value = 158

$ws.Range("D67").Value = 44295
$ws.Range("J67").Value = 2800
$ws.Range("K67").Value = 900
$ws.Range("L67").Value = 1000
$ws.Range("M67").Value = 950
$ws.Range("P67").Value = 158

$ws.Range("D68").Value = 44179
$ws.Range("J68").Value = 2200
$ws.Range("K68").Value = 800
$ws.Range("L68").Value = 1000
$ws.Range("M68").Value = 900
$ws.Range("P68").Value = 150

$ws.Range("D69").Value = 44249
$ws.Range("J69").Value = 2800
$ws.Range("K69").Value = 900
$ws.Range("L69").Value = 1000
$ws.Range("M69").Value = 950
$ws.Range("P69").Value = 158

$ws.Range("D70").Value = 44309
$ws.Range("J70").Value = 2800
$ws.Range("K70").Value = 900
$ws.Range("L70").Value = 1000
$ws.Range("M70").Value = 950
$ws.Range("P70").Value = 158

$ws.Range("D71").Value = 44384
$ws.Range("J71").Value = 3320
$ws.Range("K71").Value = 900
$ws.Range("L71").Value = 1000
$ws.Range("M71").Value = 950
$ws.Range("P71").Value = 158

$ws.Range("D72").Value = 44272
$ws.Range("J72").Value = 3100
$ws.Range("K72").Value = 800
$ws.Range("L72").Value = 1000
$ws.Range("M72").Value = 900
$ws.Range("P72").Value = 150

$ws.Range("D73").Value = 44365
$ws.Range("J73").Value = 2900
$ws.Range("K73").Value = 900
$ws.Range("L73").Value = 1000
$ws.Range("M73").Value = 950
$ws.Range("P73").Value = 158

$ws.Range("D74").Value = 44235
$ws.Range("J74").Value = 2600
$ws.Range("K74").Value = 1000
$ws.Range("L74").Value = 1100
$ws.Range("M74").Value = 1050
$ws.Range("P74").Value = 175

$ws.Range("D75").Value = 44358
$ws.Range("J75").Value = 3000
$ws.Range("K75").Value = 900
$ws.Range("L75").Value = 1000
$ws.Range("M75").Value = 950
$ws.Range("P75").Value = 158

$ws.Range("D76").Value = 44218
$ws.Range("J76").Value = 2600
$ws.Range("K76").Value = 1000
$ws.Range("L76").Value = 1100
$ws.Range("M76").Value = 1050
$ws.Range("P76").Value = 175

$ws.Range("D77").Value = 44433
$ws.Range("J77").Value = 3400
$ws.Range("K77").Value = 900
$ws.Range("L77").Value = 1000
$ws.Range("M77").Value = 950
$ws.Range("P77").Value = 158

$ws.Range("D78").Value = 44274
$ws.Range("J78").Value = 2700
$ws.Range("K78").Value = 900
$ws.Range("L78").Value = 1000
$ws.Range("M78").Value = 950
$ws.Range("P78").Value = 158

$ws.Range("D79").Value = 44321
$ws.Range("J79").Value = 3200
$ws.Range("K79").Value = 900
$ws.Range("L79").Value = 1000
$ws.Range("M79").Value = 950
$ws.Range("P79").Value = 158

$ws.Range("D80").Value = 44291
$ws.Range("J80").Value = 2600
$ws.Range("K80").Value = 900
$ws.Range("L80").Value = 1000
$ws.Range("M80").Value = 950
$ws.Range("P80").Value = 158

$ws.Range("D81").Value = 44335
$ws.Range("J81").Value = 3260
$ws.Range("K81").Value = 900
$ws.Range("L81").Value = 1000
$ws.Range("M81").Value = 950
$ws.Range("P81").Value = 158

$ws.Range("D82").Value = 44426
$ws.Range("J82").Value = 3400
$ws.Range("K82").Value = 900
$ws.Range("L82").Value = 1000
$ws.Range("M82").Value = 950
$ws.Range("P82").Value = 158

$ws.Range("D83").Value = 44421
$ws.Range("J83").Value = 3100
$ws.Range("K83").Value = 900
$ws.Range("L83").Value = 1000
$ws.Range("M83").Value = 950
$ws.Range("P83").Value = 158

$ws.Range("D84").Value = 44398
$ws.Range("J84").Value = 3360
$ws.Range("K84").Value = 900
$ws.Range("L84").Value = 1000
$ws.Range("M84").Value = 950
$ws.Range("P84").Value = 158

$ws.Range("D85").Value = 44263
$ws.Range("J85").Value = 2600
$ws.Range("K85").Value = 900
$ws.Range("L85").Value = 1000
$ws.Range("M85").Value = 950
$ws.Range("P85").Value = 158

$ws.Range("D86").Value = 44316
$ws.Range("J86").Value = 3000
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 1000
$ws.Range("M86").Value = 950
$ws.Range("P86").Value = 158

$ws.Range("D87").Value = 44253
$ws.Range("J87").Value = 2800
$ws.Range("K87").Value = 900
$ws.Range("L87").Value = 1000
$ws.Range("M87").Value = 950
$ws.Range("P87").Value = 158

$ws.Range("D88").Value = 44414
$ws.Range("J88").Value = 3200
$ws.Range("K88").Value = 900
$ws.Range("L88").Value = 1000
$ws.Range("M88").Value = 950
$ws.Range("P88").Value = 158

$ws.Range("D89").Value = 44244
$ws.Range("J89").Value = 2800
$ws.Range("K89").Value = 900
$ws.Range("L89").Value = 1000
$ws.Range("M89").Value = 950
$ws.Range("P89").Value = 158

$ws.Range("D90").Value = 44176
$ws.Range("J90").Value = 2400
$ws.Range("K90").Value = 800
$ws.Range("L90").Value = 1000
$ws.Range("M90").Value = 900
$ws.Range("P90").Value = 150

$ws.Range("D91").Value = 44417
$ws.Range("J91").Value = 3300
$ws.Range("K91").Value = 900
$ws.Range("L91").Value = 1000
$ws.Range("M91").Value = 950
$ws.Range("P91").Value = 158

$ws.Range("D92").Value = 44237
$ws.Range("J92").Value = 2800
$ws.Range("K92").Value = 1000
$ws.Range("L92").Value = 1100
$ws.Range("M92").Value = 1050
$ws.Range("P92").Value = 175

$ws.Range("D93").Value = 44351
$ws.Range("J93").Value = 2960
$ws.Range("K93").Value = 900
$ws.Range("L93").Value = 1000
$ws.Range("M93").Value = 950
$ws.Range("P93").Value = 158

$ws.Range("D94").Value = 44342
$ws.Range("J94").Value = 3300
$ws.Range("K94").Value = 900
$ws.Range("L94").Value = 1000
$ws.Range("M94").Value = 950
$ws.Range("P94").Value = 158

$ws.Range("D95").Value = 44379
$ws.Range("J95").Value = 3000
$ws.Range("K95").Value = 900
$ws.Range("L95").Value = 1000
$ws.Range("M95").Value = 950
$ws.Range("P95").Value = 158

$ws.Range("D96").Value = 44302
$ws.Range("J96").Value = 2800
$ws.Range("K96").Value = 900
$ws.Range("L96").Value = 1000
$ws.Range("M96").Value = 950
$ws.Range("P96").Value = 158

$ws.Range("D97").Value = 44391
$ws.Range("J97").Value = 3360
$ws.Range("K97").Value = 900
$ws.Range("L97").Value = 1000
$ws.Range("M97").Value = 950
$ws.Range("P97").Value = 158

$ws.Range("D98").Value = 44438
$ws.Range("J98").Value = 3100
$ws.Range("K98").Value = 900
$ws.Range("L98").Value = 1000
$ws.Range("M98").Value = 950
$ws.Range("P98").Value = 158

$ws.Range("D99").Value = 44251
$ws.Range("J99").Value = 3000
$ws.Range("K99").Value = 900
$ws.Range("L99").Value = 1000
$ws.Range("M99").Value = 950
$ws.Range("P99").Value = 158

$ws.Range("D100").Value = 44428
$ws.Range("J100").Value = 3120
$ws.Range("K100").Value = 900
$ws.Range("L100").Value = 1000
$ws.Range("M100").Value = 950
$ws.Range("P100").Value = 158

$ws.Range("D101").Value = 44442
$ws.Range("J101").Value = 3000
$ws.Range("K101").Value = 900
$ws.Range("L101").Value = 1000
$ws.Range("M101").Value = 950
$ws.Range("P101").Value = 158

$ws.Range("D102").Value = 44435
$ws.Range("J102").Value = 9720
$ws.Range("K102").Value = 900
$ws.Range("L102").Value = 1000
$ws.Range("M102").Value = 950
$ws.Range("P102").Value = 158

$ws.Range("D103").Value = 44319
$ws.Range("J103").Value = 2600
$ws.Range("K103").Value = 900
$ws.Range("L103").Value = 1000
$ws.Range("M103").Value = 950
$ws.Range("P103").Value = 158

$ws.Range("D104").Value = 44279
$ws.Range("J104").Value = 3000
$ws.Range("K104").Value = 800
$ws.Range("L104").Value = 1000
$ws.Range("M104").Value = 900
$ws.Range("P104").Value = 150

$ws.Range("D105").Value = 44412
$ws.Range("J105").Value = 3400
$ws.Range("K105").Value = 900
$ws.Range("L105").Value = 1000
$ws.Range("M105").Value = 950
$ws.Range("P105").Value = 158

$ws.Range("D106").Value = 44223
$ws.Range("J106").Value = 2800
$ws.Range("K106").Value = 1000
$ws.Range("L106").Value = 1100
$ws.Range("M106").Value = 1050
$ws.Range("P106").Value = 175

$ws.Range("D107").Value = 44314
$ws.Range("J107").Value = 3200
$ws.Range("K107").Value = 900
$ws.Range("L107").Value = 1000
$ws.Range("M107").Value = 950
$ws.Range("P107").Value = 158

$ws.Range("D108").Value = 44167
$ws.Range("J108").Value = 2700
$ws.Range("K108").Value = 800
$ws.Range("L108").Value = 1000
$ws.Range("M108").Value = 900
$ws.Range("P108").Value = 150

# New row 109: static descriptive fields copied from row 108 (identical across the dataset),
# numeric/text fields set explicitly to the shifted values.
$ws.Range("A109").Value = $ws.Range("A108").Value()
$ws.Range("B109").Value = $ws.Range("B108").Value()
$ws.Range("C109").Value = $ws.Range("C108").Value()
$ws.Range("D109").Value = 44400
$ws.Range("D109").NumberFormat = $ws.Range("D108").NumberFormat
$ws.Range("E109").Value = $ws.Range("E108").Value()
$ws.Range("F109").Value = $ws.Range("F108").Value()
$ws.Range("G109").Value = $ws.Range("G108").Value()
$ws.Range("H109").Value = $ws.Range("H108").Value()
$ws.Range("I109").Value = $ws.Range("I108").Value()
$ws.Range("J109").Value = 3100
$ws.Range("K109").Value = 900
$ws.Range("L109").Value = 1000
$ws.Range("M109").Value = 950
$ws.Range("N109").Value = $ws.Range("N108").Value()
$ws.Range("O109").Value = $ws.Range("O108").Value()
$ws.Range("P109").Value = 158
$ws.Range("Q109").Value = $ws.Range("Q108").Value()
$ws.Range("R109").Value = $ws.Range("R108").Value()
